$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.114.20"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "'1.558.97"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'0.9994"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").Value = "'291.85"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("D7").Value = "'0.3955"
$ws.Range("E7").Value = "  +4.75%  "

$ws.Range("D8").Value = "'0.3241"
$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("D9").Value = "'43.89"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").Value = "'0.07340"
$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("E11").Value = "  -3.86%  "

$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").Value = "'19.16"
$ws.Range("E13").Value = "  -5.19%  "

$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").Value = "'0.00001149"
$ws.Range("E15").Value = "  +7.00%  "

$ws.Range("D16").Value = "'6.674"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("D17").Value = "'1.555.33"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "'0.06605"

$ws.Range("E19").Value = "  -2.34%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").Value = "'6.318"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("D22").Value = "'15.85"
$ws.Range("E22").Value = "  -1.41%  "

$ws.Range("E23").Value = "  -2.44%  "

$ws.Range("D24").Value = "'22.127.07"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").Value = "'2.339"
$ws.Range("E25").Value = "  +2.15%  "

$ws.Range("D26").Value = "'2.445"
$ws.Range("E26").Value = "  -4.16%  "

$ws.Range("D27").Value = "'148.31"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = "  -2.99%  "

$ws.Range("D29").Value = "'4.874"

$ws.Range("D30").Value = "'1.729.55"
$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("D31").Value = "'119.31"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").Value = "'1.026"
$ws.Range("E32").Value = "  -4.50%  "

$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").Value = "'0.08364"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").Value = "'1.636"
$ws.Range("E35").Value = "  -13.28%  "

$ws.Range("D36").Value = "'9.096"
$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("D37").Value = "'0.06165"
$ws.Range("E37").Value = "  -1.78%  "

$ws.Range("D38").Value = "'0.02279"
$ws.Range("E38").Value = "  -3.65%  "

$ws.Range("D39").Value = "'5.167"
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").Value = "'1.219"
$ws.Range("E40").Value = "  -1.79%  "

$ws.Range("D41").Value = "'0.2070"
$ws.Range("E41").Value = "  -3.72%  "

$ws.Range("D42").Value = "'10.82"
$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("D43").Value = "'0.9987"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("D45").Value = "'13.17"
$ws.Range("E45").Value = "  -3.92%  "

$ws.Range("D46").Value = "'3.765"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "'0.5644"
$ws.Range("E47").Value = "  -4.54%  "

$ws.Range("D48").Value = "'119.16"
$ws.Range("E48").Value = "  -3.26%  "

$ws.Range("D49").Value = "'1.906"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("D50").Value = "'1.144"
$ws.Range("E50").Value = "  -2.70%  "

$ws.Range("E51").Value = "  -2.72%  "
